# "changed around stuff in meat" - lowercase all the labels/values across the
# food-group sheets, rename the sheet tabs to match (lowercase, spaces instead
# of underscores), turn the old TRUE/FALSE booleans on the "meat" sheet into
# plain lowercase text, and leave the workbook with "white meat" as the active
# tab/selection.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rewrite the data on every sheet (still addressed by their ORIGINAL names)
# ---------------------------------------------------------------------------

$meat = $wb.Worksheets.Item("Meat")
$meat.Range("A1").Value = "property"
$meat.Range("B1").Value = "value"
$meat.Range("A2").Value = "name"
$meat.Range("B2").Value = "meat"
$meat.Range("A3").Value = "healthy"
$meat.Range("B3").Value = "pass"
$meat.Range("A4").Value = "vegetarian"
$meat.Range("B4").Value = "'false"
$meat.Range("A5").Value = "lactose-free"
$meat.Range("B5").Value = "'true"
$meat.Range("A6").Value = "gluten-free"
$meat.Range("B6").Value = "'true"
$meat.Range("A7").Value = "main-category"
$meat.Range("B7").Value = "base"
$meat.Range("A8").Value = "vegan"
$meat.Range("B8").Value = "'false"

$whiteMeat = $wb.Worksheets.Item("White_Meat")
$whiteMeat.Range("A1").Value = "property"
$whiteMeat.Range("B1").Value = "value"
$whiteMeat.Range("A2").Value = "name"
$whiteMeat.Range("B2").Value = "white meat"
$whiteMeat.Range("A3").Value = "healthy"
$whiteMeat.Range("B3").Value = "yes"
$whiteMeat.Range("A4").Value = "food super group"
$whiteMeat.Range("B4").Value = "meat"
$whiteMeat.Range("A5").Value = "vegetarian substitute"
$whiteMeat.Range("B5").Value = "tofu"
$whiteMeat.Range("A6").Value = "vegan substitute"
$whiteMeat.Range("B6").Value = "tofu"

$fish = $wb.Worksheets.Item("Fish")
$fish.Range("A1").Value = "property"
$fish.Range("B1").Value = "value"
$fish.Range("A2").Value = "name"
$fish.Range("B2").Value = "fish"
$fish.Range("A3").Value = "food super group"
$fish.Range("B3").Value = "white meat"

$fowl = $wb.Worksheets.Item("Fowl")
$fowl.Range("A1").Value = "property"
$fowl.Range("B1").Value = "value"
$fowl.Range("A2").Value = "name"
$fowl.Range("B2").Value = "fowl"
$fowl.Range("A3").Value = "food super group"
$fowl.Range("B3").Value = "white meat"

$redMeat = $wb.Worksheets.Item("Red_Meat")
$redMeat.Range("A1").Value = "property"
$redMeat.Range("B1").Value = "value"
$redMeat.Range("A2").Value = "name"
$redMeat.Range("B2").Value = "red meat"
$redMeat.Range("A3").Value = "healthy"
$redMeat.Range("B3").Value = "no"
$redMeat.Range("A4").Value = "food super group"
$redMeat.Range("B4").Value = "meat"
$redMeat.Range("A5").Value = "vegetarian substitute"
$redMeat.Range("B5").Value = "plant-based red meat"
$redMeat.Range("A6").Value = "vegan substitute"
$redMeat.Range("B6").Value = "plant-based red meat"

$pork = $wb.Worksheets.Item("Pork")
$pork.Range("A1").Value = "property"
$pork.Range("B1").Value = "value"
$pork.Range("A2").Value = "name"
$pork.Range("B2").Value = "pork"
$pork.Range("A3").Value = "food super group"
$pork.Range("B3").Value = "red meat"

$beef = $wb.Worksheets.Item("Beef")
$beef.Range("A1").Value = "property"
$beef.Range("B1").Value = "value"
$beef.Range("A2").Value = "name"
$beef.Range("B2").Value = "beef"
$beef.Range("A3").Value = "food super group"
$beef.Range("B3").Value = "red meat"

$lamb = $wb.Worksheets.Item("Lamb")
$lamb.Range("A1").Value = "property"
$lamb.Range("B1").Value = "value"
$lamb.Range("A2").Value = "name"
$lamb.Range("B2").Value = "lamb"
$lamb.Range("A3").Value = "food super group"
$lamb.Range("B3").Value = "red meat"

# ---------------------------------------------------------------------------
# 2. Rename the sheet tabs themselves (lowercase, spaces instead of
#    underscores). Do this after editing the cells so the `$wb.Worksheets.Item`
#    lookups above (by original name) keep working.
# ---------------------------------------------------------------------------

$meat.Name = "meat"
$whiteMeat.Name = "white meat"
$fish.Name = "fish"
$fowl.Name = "fowl"
$redMeat.Name = "red meat"
$pork.Name = "pork"
$beef.Name = "beef"
$lamb.Name = "lamb"

# ---------------------------------------------------------------------------
# 3. Restore per-sheet selections. Selecting a range implicitly activates its
#    sheet, so the LAST selection made below ends up owning the workbook's
#    active tab -- do "white meat" last so it becomes the active sheet.
# ---------------------------------------------------------------------------

$meat.Range("C1:D8").Select()
$fish.Range("B2").Select()
$fowl.Range("A1:B3").Select()
$redMeat.Range("C1").Select()
$pork.Range("B2").Select()
$beef.Range("A1:B3").Select()
$lamb.Range("A1:B3").Select()
$whiteMeat.Range("C15").Select()
